# Commit: "changed get_involved route to get-involved"
#
# The "Overview" tab instructs contributors to upload their completed
# form by visiting a "get_involved" URL; the route was renamed to use a
# hyphen ("get-involved") instead of an underscore. Update the cell text
# accordingly.

$wb = $excel.ActiveWorkbook

$old = "To upload your completed form, follow the instructions at https://openenzymedb.platform.moleculemaker.org/about/get_involved."
$new = "To upload your completed form, follow the instructions at https://openenzymedb.platform.moleculemaker.org/about/get-involved."

$found = $false

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    if ($used -eq $null) { continue }
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value2 -eq $old) {
                $cell.Value = $new
                $found = $true
            }
        }
    }
}

if (-not $found) {
    throw "Could not locate the cell containing the get_involved instructions."
}
